# chore: update Sheets via scheduled runner
# Refreshes the computed market-price / profit columns (H:N) on each of the
# per-job sheets with newly scraped values. Columns:
#   H currentAveragePrice      I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ              L LevePriceHQ             M LeveProfitNQ
#   N LeveProfitHQ
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value2 = 907.7273
$ws.Cells.Item(32, 10).Value2 = 830.8333
$ws.Cells.Item(32, 12).Value2 = 830.8333
$ws.Cells.Item(32, 14).Value2 = -1482.8333
$ws.Cells.Item(34, 8).Value2 = 3274.0908
$ws.Cells.Item(34, 9).Value2 = 3274.0908
$ws.Cells.Item(34, 11).Value2 = 3274.0908
$ws.Cells.Item(34, 13).Value2 = -3071.0908
$ws.Cells.Item(36, 8).Value2 = 3274.0908
$ws.Cells.Item(36, 9).Value2 = 3274.0908
$ws.Cells.Item(36, 11).Value2 = 3274.0908
$ws.Cells.Item(36, 13).Value2 = -2559.0908
$ws.Cells.Item(99, 8).Value2 = 433.33334
$ws.Cells.Item(99, 10).Value2 = 500
$ws.Cells.Item(99, 12).Value2 = 1500
$ws.Cells.Item(99, 14).Value2 = -4496
$ws.Cells.Item(107, 8).Value2 = 621.6667
$ws.Cells.Item(107, 9).Value2 = 687
$ws.Cells.Item(107, 11).Value2 = 687
$ws.Cells.Item(107, 13).Value2 = 1233
$ws.Cells.Item(129, 8).Value2 = 1805.3636
$ws.Cells.Item(129, 9).Value2 = 551.4286
$ws.Cells.Item(129, 11).Value2 = 1654.2858
$ws.Cells.Item(129, 13).Value2 = 3345.7142
$ws.Cells.Item(137, 8).Value2 = 1889.2307
$ws.Cells.Item(137, 9).Value2 = 1760.8334
$ws.Cells.Item(137, 10).Value2 = 1999.2858
$ws.Cells.Item(137, 11).Value2 = 5282.5002
$ws.Cells.Item(137, 12).Value2 = 5997.857400000001
$ws.Cells.Item(137, 13).Value2 = -2732.5002
$ws.Cells.Item(137, 14).Value2 = -11097.8574
$ws.Cells.Item(138, 8).Value2 = 4389.244
$ws.Cells.Item(138, 9).Value2 = 3675.6086
$ws.Cells.Item(138, 10).Value2 = 5301.1113
$ws.Cells.Item(138, 11).Value2 = 11026.8258
$ws.Cells.Item(138, 12).Value2 = 15903.3339
$ws.Cells.Item(138, 13).Value2 = -5886.825800000001
$ws.Cells.Item(138, 14).Value2 = -26183.3339

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 5986.273
$ws.Cells.Item(32, 9).Value2 = 4323.077
$ws.Cells.Item(32, 11).Value2 = 4323.077
$ws.Cells.Item(32, 13).Value2 = -4036.077
$ws.Cells.Item(43, 8).Value2 = 29999.334
$ws.Cells.Item(43, 9).Value2 = 30000
$ws.Cells.Item(43, 11).Value2 = 30000
$ws.Cells.Item(43, 13).Value2 = -29687
$ws.Cells.Item(45, 8).Value2 = 3761.8
$ws.Cells.Item(45, 9).Value2 = 3639.75
$ws.Cells.Item(45, 10).Value2 = 4250
$ws.Cells.Item(45, 11).Value2 = 3639.75
$ws.Cells.Item(45, 12).Value2 = 4250
$ws.Cells.Item(45, 13).Value2 = -3262.75
$ws.Cells.Item(45, 14).Value2 = -5004
$ws.Cells.Item(61, 8).Value2 = 1889.75
$ws.Cells.Item(61, 9).Value2 = 1889.75
$ws.Cells.Item(61, 11).Value2 = 1889.75
$ws.Cells.Item(61, 13).Value2 = -1677.75
$ws.Cells.Item(122, 8).Value2 = 2182.75
$ws.Cells.Item(122, 10).Value2 = 3831.6667
$ws.Cells.Item(122, 12).Value2 = 11495.0001
$ws.Cells.Item(122, 14).Value2 = -16395.0001
$ws.Cells.Item(132, 8).Value2 = 727
$ws.Cells.Item(132, 9).Value2 = 514.8333
$ws.Cells.Item(132, 11).Value2 = 1544.4999
$ws.Cells.Item(132, 13).Value2 = 985.5001
$ws.Cells.Item(136, 8).Value2 = 1889.75
$ws.Cells.Item(136, 9).Value2 = 1889.75
$ws.Cells.Item(136, 11).Value2 = 5669.25
$ws.Cells.Item(136, 13).Value2 = -3119.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value2 = 3106.3076
$ws.Cells.Item(105, 9).Value2 = 2853
$ws.Cells.Item(105, 11).Value2 = 2853
$ws.Cells.Item(105, 13).Value2 = -1106
$ws.Cells.Item(107, 8).Value2 = 8955
$ws.Cells.Item(107, 9).Value2 = 8955
$ws.Cells.Item(107, 11).Value2 = 8955
$ws.Cells.Item(107, 13).Value2 = -7035

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 3725.3914
$ws.Cells.Item(31, 9).Value2 = 3189.3333
$ws.Cells.Item(31, 11).Value2 = 3189.3333
$ws.Cells.Item(31, 13).Value2 = -2894.3333
$ws.Cells.Item(34, 8).Value2 = 3725.3914
$ws.Cells.Item(34, 9).Value2 = 3189.3333
$ws.Cells.Item(34, 11).Value2 = 3189.3333
$ws.Cells.Item(34, 13).Value2 = -2987.3333
$ws.Cells.Item(58, 8).Value2 = 1935.4762
$ws.Cells.Item(58, 9).Value2 = 1103.6177
$ws.Cells.Item(58, 11).Value2 = 1103.6177
$ws.Cells.Item(58, 13).Value2 = -900.6177
$ws.Cells.Item(62, 8).Value2 = 38595.273
$ws.Cells.Item(62, 9).Value2 = 2349.8
$ws.Cells.Item(62, 11).Value2 = 2349.8
$ws.Cells.Item(62, 13).Value2 = -1725.8
$ws.Cells.Item(65, 8).Value2 = 38595.273
$ws.Cells.Item(65, 9).Value2 = 2349.8
$ws.Cells.Item(65, 11).Value2 = 11749
$ws.Cells.Item(65, 13).Value2 = -8629
$ws.Cells.Item(122, 8).Value2 = 2427.5
$ws.Cells.Item(122, 9).Value2 = 2581.842
$ws.Cells.Item(122, 11).Value2 = 7745.526
$ws.Cells.Item(122, 13).Value2 = -5295.526
$ws.Cells.Item(129, 8).Value2 = 0
$ws.Cells.Item(129, 10).Value2 = 0
$ws.Cells.Item(129, 12).Value2 = 0
$ws.Cells.Item(129, 14).ClearContents() | Out-Null
$ws.Cells.Item(132, 8).Value2 = 1544.2858
$ws.Cells.Item(132, 9).Value2 = 1321.2307
$ws.Cells.Item(132, 11).Value2 = 3963.6921
$ws.Cells.Item(132, 13).Value2 = -1433.6921
$ws.Cells.Item(134, 8).Value2 = 1767.9744
$ws.Cells.Item(134, 9).Value2 = 1573
$ws.Cells.Item(134, 11).Value2 = 4719
$ws.Cells.Item(134, 13).Value2 = -2184
$ws.Cells.Item(136, 8).Value2 = 1935.4762
$ws.Cells.Item(136, 9).Value2 = 1103.6177
$ws.Cells.Item(136, 11).Value2 = 3310.8531
$ws.Cells.Item(136, 13).Value2 = -760.8531000000003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value2 = 35580.633
$ws.Cells.Item(122, 9).Value2 = 1737.0435
$ws.Cells.Item(122, 10).Value2 = 146781
$ws.Cells.Item(122, 11).Value2 = 5211.1305
$ws.Cells.Item(122, 12).Value2 = 440343
$ws.Cells.Item(122, 13).Value2 = -2761.1305
$ws.Cells.Item(122, 14).Value2 = -445243
$ws.Cells.Item(132, 8).Value2 = 2115.138
$ws.Cells.Item(132, 9).Value2 = 1892.6086
$ws.Cells.Item(132, 11).Value2 = 5677.825800000001
$ws.Cells.Item(132, 13).Value2 = -3147.825800000001
$ws.Cells.Item(141, 8).Value2 = 78995
$ws.Cells.Item(141, 10).Value2 = 78995
$ws.Cells.Item(141, 12).Value2 = 78995
$ws.Cells.Item(141, 14).Value2 = -89355

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 3999
$ws.Cells.Item(46, 9).Value2 = 2928.5715
$ws.Cells.Item(46, 10).Value2 = 5872.25
$ws.Cells.Item(46, 11).Value2 = 2928.5715
$ws.Cells.Item(46, 12).Value2 = 5872.25
$ws.Cells.Item(46, 13).Value2 = -2740.5715
$ws.Cells.Item(46, 14).Value2 = -6248.25
$ws.Cells.Item(93, 8).Value2 = 2186.7144
$ws.Cells.Item(93, 9).Value2 = 1802
$ws.Cells.Item(93, 10).Value2 = 2475.25
$ws.Cells.Item(93, 11).Value2 = 1802
$ws.Cells.Item(93, 12).Value2 = 2475.25
$ws.Cells.Item(93, 13).Value2 = -554
$ws.Cells.Item(93, 14).Value2 = -4971.25
$ws.Cells.Item(100, 8).Value2 = 4162.4443
$ws.Cells.Item(100, 9).Value2 = 4162.4443
$ws.Cells.Item(100, 11).Value2 = 4162.4443
$ws.Cells.Item(100, 13).Value2 = -3621.4443

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value2 = 1777.5
$ws.Cells.Item(100, 9).Value2 = 1888.5714
$ws.Cells.Item(100, 10).Value2 = 1000
$ws.Cells.Item(100, 11).Value2 = 3777.1428
$ws.Cells.Item(100, 12).Value2 = 2000
$ws.Cells.Item(100, 13).Value2 = -3236.1428
$ws.Cells.Item(100, 14).Value2 = -3082
$ws.Cells.Item(132, 8).Value2 = 5277.3335
$ws.Cells.Item(132, 9).Value2 = 3297.375
$ws.Cells.Item(132, 11).Value2 = 9892.125
$ws.Cells.Item(132, 13).Value2 = -7362.125
